{"js": "// Add <s>/<sub>/<sup> style formatting (strikethrough, subscript, superscript)\n// to specific words inside the long \"Text\" paragraph, matching the target\n// OOXML run split:\n//   ... dui [egestas]{strike,sup}[, ]{strike}[volutpat]{strike,sub} nisi ...\n//   ... augue elit [aliquam ]{sup}[mauris]{sub}, vel mollis ... ut [ipsum]{strike}.\n\n// Anchor on a long, unique phrase from the target paragraph so we never\n// confuse it with similar words/phrases that occur elsewhere in the document.\nconst anchorHits = context.document.body.search(\n  \"Praesent ornare fermentum turpis\",\n  { matchCase: true }\n);\nanchorHits.load(\"items\");\nawait context.sync();\n\nconst paragraph = anchorHits.items[0].paragraphs.getFirst();\n\n// Each target word shows up more than once inside the paragraph (e.g.\n// \"egestas\" twice, \"aliquam\" three times, \"ipsum\" four times), so a plain\n// paragraph.search(\"egestas\") would be ambiguous. Instead we first search for\n// a short, unique surrounding phrase, then search *inside* that hit for the\n// exact word/substring we want to format. That two-step (phrase -> word)\n// search pins down exactly the right occurrence.\nconst egestasPhrase = paragraph.search(\"dui egestas\", { matchCase: true });\nconst commaPhrase = paragraph.search(\"egestas, volutpat\", { matchCase: true });\nconst volutpatPhrase = paragraph.search(\", volutpat\", { matchCase: true });\nconst aliquamPhrase = paragraph.search(\"elit aliquam mauris\", { matchCase: true });\nconst maurisPhrase = paragraph.search(\"aliquam mauris\", { matchCase: true });\nconst ipsumPhrase = paragraph.search(\"ut ipsum.\", { matchCase: true });\nawait context.sync();\n\nconst egestasWord = egestasPhrase.items[0].search(\"egestas\", { matchCase: true });\nconst commaWord = commaPhrase.items[0].search(\", \", { matchCase: true });\nconst volutpatWord = volutpatPhrase.items[0].search(\"volutpat\", { matchCase: true });\nconst aliquamWord = aliquamPhrase.items[0].search(\"aliquam \", { matchCase: true });\nconst maurisWord = maurisPhrase.items[0].search(\"mauris\", { matchCase: true });\nconst ipsumWord = ipsumPhrase.items[0].search(\"ipsum\", { matchCase: true });\nawait context.sync();\n\n// 1. \"egestas\" -> <s><sup>\negestasWord.items[0].font.strikeThrough = true;\negestasWord.items[0].font.superscript = true;\n\n// 2. \", \" between egestas and volutpat -> <s>\ncommaWord.items[0].font.strikeThrough = true;\n\n// 3. \"volutpat\" -> <s><sub>\nvolutpatWord.items[0].font.strikeThrough = true;\nvolutpatWord.items[0].font.subscript = true;\n\n// 4. \"aliquam \" -> <sup>\naliquamWord.items[0].font.superscript = true;\n\n// 5. \"mauris\" -> <sub>\nmaurisWord.items[0].font.subscript = true;\n\n// 6. \"ipsum\" (the last one, right before the final period) -> <s>\nipsumWord.items[0].font.strikeThrough = true;\n\nawait context.sync();\n", "ps1": "# Add <s>/<sub>/<sup> style formatting (strikethrough, subscript, superscript)\n# to specific words inside the long \"Text\" paragraph, matching the target\n# OOXML run split:\n#   ... dui [egestas]{strike,sup}[, ]{strike}[volutpat]{strike,sub} nisi ...\n#   ... augue elit [aliquam ]{sup}[mauris]{sub}, vel mollis ... ut [ipsum]{strike}.\n\n$d = $word.ActiveDocument\n\n# Anchor on a long, unique phrase from the target paragraph so we never\n# confuse it with similar words/phrases that occur elsewhere in the document.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Praesent ornare fermentum turpis*\") {\n        $target = $p.Range\n        break\n    }\n}\n\n# Helper: search for $phrase inside $scopeRange and return a fresh Range\n# collapsed onto the hit (a Duplicate is used so the original scope range is\n# left untouched for later reuse).\nfunction Find-SubRange($scopeRange, $phrase) {\n    $r = $scopeRange.Duplicate\n    $r.Find.ClearFormatting()\n    $null = $r.Find.Execute($phrase, $true)\n    return $r\n}\n\n# Each target word shows up more than once inside the paragraph (e.g.\n# \"egestas\" twice, \"aliquam\" three times, \"ipsum\" four times), so searching\n# for the bare word directly would be ambiguous. Instead we first find a\n# short, unique surrounding phrase, then search *inside* that hit for the\n# exact word/substring we want to format. That two-step (phrase -> word)\n# search pins down exactly the right occurrence.\n\n# 1. \"egestas\" -> <s><sup>\n$egestasPhrase = Find-SubRange $target \"dui egestas\"\n$egestasWord = Find-SubRange $egestasPhrase \"egestas\"\n$egestasWord.Font.StrikeThrough = $true\n$egestasWord.Font.Superscript = $true\n\n# 2. \", \" between egestas and volutpat -> <s>\n$commaPhrase = Find-SubRange $target \"egestas, volutpat\"\n$commaWord = Find-SubRange $commaPhrase \", \"\n$commaWord.Font.StrikeThrough = $true\n\n# 3. \"volutpat\" -> <s><sub>\n$volutpatPhrase = Find-SubRange $target \", volutpat\"\n$volutpatWord = Find-SubRange $volutpatPhrase \"volutpat\"\n$volutpatWord.Font.StrikeThrough = $true\n$volutpatWord.Font.Subscript = $true\n\n# 4. \"aliquam \" -> <sup>\n$aliquamPhrase = Find-SubRange $target \"elit aliquam mauris\"\n$aliquamWord = Find-SubRange $aliquamPhrase \"aliquam \"\n$aliquamWord.Font.Superscript = $true\n\n# 5. \"mauris\" -> <sub>\n$maurisPhrase = Find-SubRange $target \"aliquam mauris\"\n$maurisWord = Find-SubRange $maurisPhrase \"mauris\"\n$maurisWord.Font.Subscript = $true\n\n# 6. \"ipsum\" (the last one, right before the final period) -> <s>\n$ipsumPhrase = Find-SubRange $target \"ut ipsum.\"\n$ipsumWord = Find-SubRange $ipsumPhrase \"ipsum\"\n$ipsumWord.Font.StrikeThrough = $true\n"}
